# Scene 23C line-edits
#
# 1) Narration: "maybe his parent is inside?" -> "maybe his parent's inside?"
# 2) Narration: "even if his parent is inside the store" ->
#               "even if his parent'ss inside the store"   (typo kept verbatim, per the edit)
# 3) Stage direction: "?Greta (embarrassed embarrassed): I... uh..." ->
#                     "?Greta (neutral embarrassed): I... uh..."

$d = $word.ActiveDocument

$ok1 = $d.Content.Find.Execute(
    "maybe his parent is inside?",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "maybe his parent’s inside?",
    2
)

$ok2 = $d.Content.Find.Execute(
    "even if his parent is inside the store",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "even if his parent’ss inside the store",
    2
)

$ok3 = $d.Content.Find.Execute(
    "?Greta (embarrassed embarrassed): I… uh…",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "?Greta (neutral embarrassed): I… uh…",
    2
)

"parent's-inside replace: $ok1"
"parent'ss-inside replace: $ok2"
"Greta neutral-embarrassed replace: $ok3"
